$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" column header
$ws.Cells.Item(1, 8).Value = "Label"

# Copy style from G1 (bold/bordered header style) to H1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

# Fill in the new "Label" column values (0 = Control, 1 = MDD) for both blocks (rows 2-11, 12-21)
$labelValues = @(0,0,0,0,0,1,1,1,1,1,0,0,0,0,0,1,1,1,1,1)
for ($i = 0; $i -lt $labelValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $labelValues[$i]
}

# Update refit values for the individually-refit first block (rows 4-11)
$ws.Range("D4").Value = 0.3932534943075598
$ws.Range("E4").Value = 0.3932534943075598

$ws.Range("D5").Value = 0.4492191791131242
$ws.Range("E5").Value = 0.4492191791131242

$ws.Range("D6").Value = 0.4569337603861072
$ws.Range("E6").Value = 0.4569337603861072

$ws.Range("D9").Value = 0.5540508035896428
$ws.Range("E9").Value = 0.4459491964103572

$ws.Range("D11").Value = 0.3479577628909745
$ws.Range("E11").Value = 0.6520422371090255
$ws.Range("F11").Value = 0.6623891592025757
